$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells for numeric-looking price strings so Excel
# does not silently coerce them to Double (losing trailing zeros / exact text).

# Row 2
$ws.Range('D2').Value = '66.730.59'
$ws.Range('E2').Value = '  +0.90%  '

# Row 3
$ws.Range('D3').Value = '3.356.54'
$ws.Range('E3').Value = '  +1.35%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '190.60'
$ws.Range('E5').Value = '  +4.86%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '562.76'
$ws.Range('E6').Value = '  +0.50%  '

# Row 7
$ws.Range('E7').Value = '  -0.04%  '

# Row 8
$ws.Range('D8').Value = '3.351.17'
$ws.Range('E8').Value = '  +1.36%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.584'
$ws.Range('E9').Value = '  -1.23%  '

# Row 10
$ws.Range('E10').Value = '  -2.87%  '

# Row 11
$ws.Range('E11').Value = '  -0.29%  '

# Row 12
$ws.Range('E12').Value = '  -1.27%  '

# Row 13
$ws.Range('E13').Value = '  +1.95%  '

# Row 14
$ws.Range('E14').Value = '  +1.81%  '

# Row 15
$ws.Range('D15').Value = '3.890.04'
$ws.Range('E15').Value = '  +1.44%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '605.48'
$ws.Range('E16').Value = '  -5.17%  '

# Row 17
$ws.Range('D17').Value = '66.776.01'
$ws.Range('E17').Value = '  +1.10%  '

# Row 18
$ws.Range('E18').Value = '  +0.12%  '

# Row 19
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.358.43'
$ws.Range('E19').Value = '  +1.71%  '

# Row 20
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.118'
$ws.Range('E20').Value = '  +1.19%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.11'
$ws.Range('E21').Value = '  -3.15%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.908'
$ws.Range('E22').Value = '  +0.15%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.45'
$ws.Range('E23').Value = '  +4.52%  '

# Row 24
$ws.Range('E24').Value = '  +0.02%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '101.05'
$ws.Range('E25').Value = '  -5.87%  '

# Row 26
$ws.Range('E26').Value = '  +0.17%  '

# Row 27
$ws.Range('E27').Value = '  +1.16%  '

# Row 28
$ws.Range('E28').Value = '  +2.83%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.69'
$ws.Range('E29').Value = '  +1.57%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.76'
$ws.Range('E30').Value = '  -0.33%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.89'
$ws.Range('E31').Value = '  +0.35%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.79'
$ws.Range('E32').Value = '  +6.92%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.00'
$ws.Range('E33').Value = '  +0.00%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '586.50'
$ws.Range('E34').Value = '  +5.28%  '

# Row 35
$ws.Range('E35').Value = '  +0.04%  '

# Row 36
$ws.Range('E36').Value = '  +0.07%  '

# Row 37
$ws.Range('D37').Value = '3.736.02'
$ws.Range('E37').Value = '  -0.05%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '57.21'
$ws.Range('E38').Value = '  -0.75%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.00%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.59'
$ws.Range('E40').Value = '  +6.95%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '34.06'
$ws.Range('E41').Value = '  +5.70%  '

# Row 42
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.130'
$ws.Range('E42').Value = '  +1.43%  '

# Row 43
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0714'
$ws.Range('E43').Value = '  -0.74%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.27'
$ws.Range('E44').Value = '  -7.20%  '

# Row 45
$ws.Range('E45').Value = '  -1.00%  '

# Row 46
$ws.Range('E46').Value = '  +0.42%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.39'
$ws.Range('E47').Value = '  +4.74%  '

# Row 48
$ws.Range('E48').Value = '  +2.30%  '

# Row 49
$ws.Range('E49').Value = '  +0.16%  '

# Row 50
$ws.Range('E50').Value = '  -0.85%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.31%  '
